# The author corrected the body-text bullet that credits "Tomas Muller" on
# the slide with sldId 259 (shape id=3, "Text Placeholder 2"): the words
# "by Tomas " were retyped as "by z` " (mid-edit slip), which splits the
# single run into three runs:
#   "Taking the ideas proposed " | "by z` " | "Muller, trying ... algorithm. "
$p = $ppt.ActivePresentation

# Find the target slide by its persistent SlideID (259) instead of assuming
# a fixed position in the deck.
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 259) {
        $targetSlide = $candidate
        break
    }
}
if ($targetSlide -eq $null) {
    $targetSlide = $p.Slides.Item(2)
}

# Find the target shape by its persistent shape Id (3) instead of assuming
# a fixed shape index.
$targetShape = $null
for ($i = 1; $i -le $targetSlide.Shapes.Count; $i++) {
    $candidate = $targetSlide.Shapes.Item($i)
    if ($candidate.Id -eq 3) {
        $targetShape = $candidate
        break
    }
}
if ($targetShape -eq $null) {
    $targetShape = $targetSlide.Shapes.Item(2)
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text
$target = "by Tomas "
$idx = $fullText.IndexOf($target)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $target.Length)
    $sub.Text = "by z`` "
}
